$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SFBR60_RE")

$flowType = @{
    4 = "rise"
    5 = "rise"
    6 = "peak"
    7 = "peak"
    8 = "fall"
    9 = "fall"
    10 = ""
    11 = ""
    12 = ""
    13 = "base"
    14 = "rise"
    15 = "rise"
    16 = "rise"
    17 = "rise"
    18 = "rise"
    19 = "rise"
    20 = ""
    21 = ""
    22 = ""
    23 = "peak"
    24 = "peak"
    25 = "peak"
    26 = "peak"
    27 = "fall"
    28 = "fall"
    29 = "fall"
    30 = "fall"
    31 = ""
    32 = ""
    33 = "fall"
    34 = "fall"
    35 = "base"
    36 = "rise"
    37 = "rise"
    38 = "rise"
    39 = "rise"
    40 = ""
    41 = ""
    42 = ""
    43 = "peak"
    44 = "peak"
    45 = "peak"
    46 = "fall"
    47 = "fall"
    48 = "fall"
    49 = "fall"
    50 = "fall"
    51 = ""
    52 = ""
    53 = "base"
    54 = "rise"
    55 = "rise"
    56 = "rise"
    57 = "rise"
    58 = "rise"
    59 = "rise"
    60 = "rise"
    61 = "peak"
    62 = "peak"
    63 = "fall"
    64 = "fall"
    65 = "fall"
    66 = "fall"
    67 = "base"
    68 = "base"
    69 = ""
    70 = ""
    71 = ""
    72 = ""
    73 = ""
    74 = "rise"
    75 = "peak"
    76 = "peak"
    77 = "peak"
    78 = "peak"
    79 = "fall"
    80 = "fall"
    81 = "fall"
    82 = ""
    83 = ""
    84 = "base"
    85 = "base"
    86 = "base"
    87 = "base"
    88 = "base"
    89 = "base"
    90 = "base"
    91 = "base"
    92 = "base"
    93 = "base"
    94 = ""
    95 = ""
    96 = ""
    97 = ""
    98 = ""
    99 = ""
    100 = ""
    101 = ""
    102 = ""
    103 = ""
    104 = ""
    105 = ""
    106 = ""
    107 = ""
    108 = ""
    109 = ""
    110 = ""
    111 = ""
    112 = ""
    113 = ""
    114 = ""
    115 = ""
    116 = ""
    117 = ""
    118 = ""
    119 = ""
    120 = ""
    121 = ""
    122 = ""
    123 = ""
    124 = ""
    125 = ""
    126 = ""
    127 = ""
    128 = ""
    129 = ""
    130 = ""
    131 = ""
    132 = ""
    133 = ""
    134 = ""
    135 = ""
    136 = ""
    137 = ""
    138 = ""
    139 = ""
    140 = ""
    141 = ""
    142 = ""
    143 = ""
    144 = "base"
    145 = "base"
    146 = "base"
    147 = "base"
}

# Populate data rows first (so shared-string insertion order matches
# rise/peak/fall/base before the "Flow_Type" header text)
foreach ($row in ($flowType.Keys | Sort-Object {[int]$_})) {
    $val = $flowType[$row]
    if ($val -ne "") {
        $ws.Cells.Item([int]$row, 8).Value2 = $val
    }
}

# Header goes in last so it becomes the final new shared string
$ws.Cells.Item(1, 8).Value2 = "Flow_Type"
